# Rapise 7.4 update: TC_CreateNewContact/Main.rvl
# - "Contacts" click action becomes a left-click (DoLClick) in both the
#   RVL and Cleanup sheets.
# - In RVL, a new "Last_Name" click step is inserted before the existing
#   "Last_Name" text-entry step, the old "CloseItem" step is dropped, and
#   a "DoSleep 3000" step is appended instead.
# - In Cleanup, App_Launcher1/View_All1 are renamed (drop the trailing
#   "1") and the final LIGHTNING-ICON click becomes a ClearFilter click.

$wb = $excel.ActiveWorkbook

$rvl = $wb.Worksheets.Item("RVL")
$cleanup = $wb.Worksheets.Item("Cleanup")

# --- RVL sheet ---------------------------------------------------------

# Contacts: DoClick -> DoLClick
$rvl.Cells.Item(16, 4).Value = "DoLClick"

# Shift rows 21-23 down to 22-24, inserting a new "Last_Name" / DoClick
# row at 21, renaming the old "*Last_Name" row to "Last_Name", keeping
# "SaveEdit" as-is, and replacing "CloseItem" with a "DoSleep" step.

# Row 21 (new): Action | Last_Name | DoClick
$rvl.Cells.Item(21, 1).ClearContents()
$rvl.Cells.Item(21, 2).Value = "Action"
$rvl.Cells.Item(21, 3).Value = "Last_Name"
$rvl.Cells.Item(21, 4).Value = "DoClick"
$rvl.Cells.Item(21, 5).ClearContents()
$rvl.Cells.Item(21, 6).ClearContents()
$rvl.Cells.Item(21, 7).ClearContents()

# Row 22 (was row 21): Action | Last_Name | DoSetText | txt | Data | LastName
$rvl.Cells.Item(22, 2).Value = "Action"
$rvl.Cells.Item(22, 3).Value = "Last_Name"
$rvl.Cells.Item(22, 4).Value = "DoSetText"
$rvl.Cells.Item(22, 5).Value = "txt"
$rvl.Cells.Item(22, 6).Value = "Data"
$rvl.Cells.Item(22, 7).Value = "LastName"

# Row 23 (was row 22): Action | SaveEdit | DoClick
$rvl.Cells.Item(23, 2).Value = "Action"
$rvl.Cells.Item(23, 3).Value = "SaveEdit"
$rvl.Cells.Item(23, 4).Value = "DoClick"
$rvl.Cells.Item(23, 5).ClearContents()
$rvl.Cells.Item(23, 6).ClearContents()
$rvl.Cells.Item(23, 7).ClearContents()

# Row 24 (was blank, "CloseItem" step removed): Action | Global | DoSleep | millis | number | 3000
$rvl.Cells.Item(24, 2).Value = "Action"
$rvl.Cells.Item(24, 3).Value = "Global"
$rvl.Cells.Item(24, 4).Value = "DoSleep"
$rvl.Cells.Item(24, 5).Value = "millis"
$rvl.Cells.Item(24, 6).Value = "number"
$rvl.Cells.Item(24, 7).NumberFormat = "@"
$rvl.Cells.Item(24, 7).Value = "3000"
$rvl.Cells.Item(24, 7).Style = "Normal"

# --- Cleanup sheet -------------------------------------------------------

$cleanup.Cells.Item(11, 3).Value = "App_Launcher"
$cleanup.Cells.Item(12, 3).Value = "View_All"
$cleanup.Cells.Item(13, 4).Value = "DoLClick"
$cleanup.Cells.Item(22, 3).Value = "ClearFilter"
